$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 23.40000000000022
$ws.Range("H2").Value = 0.0004705476899340422
$ws.Range("I2").Value = 0.0004705476899340422
$ws.Range("L2").Value = 48.2156695533307
$ws.Range("M2").Value = "[19.331988134786585, 77.09935097187481]"
$ws.Range("N2").Value = 0.001586487413218762
$ws.Range("O2").Value = 0.001586487413218762
$ws.Range("P2").Value = 2.081816152829272
$ws.Range("Q2").Value = "[1.452868674633116, 2.7107636310254275]"
$ws.Range("R2").Value = [double]"3.185510855097107e-08"
$ws.Range("S2").Value = [double]"3.185510855097107e-08"
$ws.Range("T2").Value = 70.47216119945074
$ws.Range("U2").Value = "[54.507361001300836, 86.43696139760064]"
$ws.Range("V2").Value = [double]"1.809019600784723e-11"
$ws.Range("W2").Value = [double]"1.809019600784723e-11"
$ws.Range("X2").Value = 15.64684684684699
$ws.Range("Y2").Value = 13.30450450450463
$ws.Range("Z2").Value = 17.98918918918936
$ws.Range("F3").Value = 23.40000000000022
$ws.Range("H3").Value = 0.0008698441710963678
$ws.Range("I3").Value = 0.0008698441710963678
$ws.Range("L3").Value = 51.13748879053259
$ws.Range("M3").Value = "[18.45151923564255, 83.82345834542262]"
$ws.Range("N3").Value = 0.00289127774284248
$ws.Range("O3").Value = 0.00289127774284248
$ws.Range("P3").Value = 1.867974010242579
$ws.Range("Q3").Value = "[1.1761317842268086, 2.5598162362583503]"
$ws.Range("R3").Value = [double]"2.112949634813788e-06"
$ws.Range("S3").Value = [double]"2.112949634813788e-06"
$ws.Range("T3").Value = 70.55929907617387
$ws.Range("U3").Value = "[52.84045278286165, 88.2781453694861]"
$ws.Range("V3").Value = [double]"3.218456612330556e-10"
$ws.Range("W3").Value = [double]"3.218456612330556e-10"
$ws.Range("X3").Value = 16.4432432432434
$ws.Range("Y3").Value = 13.8666666666668
$ws.Range("Z3").Value = 19.01981981982
$ws.Range("B4").Value = 0
$ws.Range("F4").Value = 23.40000000000022
$ws.Range("H4").Value = 0.02954131888024802
$ws.Range("I4").Value = 0.02954131888024802
$ws.Range("L4").Value = 32.50817241462209
$ws.Range("M4").Value = "[-0.2772523798346924, 65.29359720907888]"
$ws.Range("N4").Value = 0.05188640582727944
$ws.Range("O4").Value = 0.05188640582727944
$ws.Range("P4").Value = 2.119553001521042
$ws.Range("Q4").Value = "[0.786184347745194, 3.4529216552968895]"
$ws.Range("R4").Value = 0.00250831651464245
$ws.Range("S4").Value = 0.00250831651464245
$ws.Range("T4").Value = 68.72106727633715
$ws.Range("U4").Value = "[52.04372151993737, 85.39841303273693]"
$ws.Range("V4").Value = [double]"1.269153671046297e-10"
$ws.Range("W4").Value = [double]"1.269153671046297e-10"
$ws.Range("X4").Value = 15.50630630630645
$ws.Range("Y4").Value = 10.54054054054064
$ws.Range("Z4").Value = 20.47207207207226
$ws.Range("F5").Value = 23.40000000000022
$ws.Range("H5").Value = 0.001994370665514666
$ws.Range("I5").Value = 0.001994370665514666
$ws.Range("L5").Value = 50.96092381659204
$ws.Range("M5").Value = "[15.592183256059947, 86.32966437712412]"
$ws.Range("N5").Value = 0.005721876538554715
$ws.Range("O5").Value = 0.005721876538554715
$ws.Range("P5").Value = 2.283079345852042
$ws.Range("Q5").Value = "[1.628973968528041, 2.937184723176043]"
$ws.Range("R5").Value = [double]"9.19811427024797e-09"
$ws.Range("S5").Value = [double]"9.19811427024797e-09"
$ws.Range("T5").Value = 69.6206683279826
$ws.Range("U5").Value = "[50.47295207533702, 88.76838458062818]"
$ws.Range("V5").Value = [double]"3.387468749949107e-09"
$ws.Range("W5").Value = [double]"3.387468749949107e-09"
$ws.Range("X5").Value = 14.89729729729744
$ws.Range("Y5").Value = 12.46126126126138
$ws.Range("Z5").Value = 17.3333333333335
$ws.Range("F6").Value = 23.40000000000022
$ws.Range("H6").Value = [double]"1.029253295481425e-06"
$ws.Range("I6").Value = [double]"1.029253295481425e-06"
$ws.Range("L6").Value = 63.31950996972584
$ws.Range("M6").Value = "[35.05133423956423, 91.58768569988746]"
$ws.Range("N6").Value = [double]"4.585817210411847e-05"
$ws.Range("O6").Value = [double]"4.585817210411847e-05"
$ws.Range("P6").Value = 2.371131992799504
$ws.Range("Q6").Value = "[1.956026657190041, 2.786237328408966]"
$ws.Range("R6").Value = [double]"5.329070518200751e-15"
$ws.Range("S6").Value = [double]"5.329070518200751e-15"
$ws.Range("T6").Value = 66.03048803615665
$ws.Range("U6").Value = "[51.419790267881154, 80.64118580443215]"
$ws.Range("V6").Value = [double]"9.088285679581531e-12"
$ws.Range("W6").Value = [double]"9.088285679581531e-12"
$ws.Range("X6").Value = 14.56936936936951
$ws.Range("Y6").Value = 13.02342342342354
$ws.Range("Z6").Value = 16.11531531531547
$ws.Range("F7").Value = 23.40000000000022
$ws.Range("H7").Value = [double]"1.089786613794885e-05"
$ws.Range("I7").Value = [double]"1.089786613794885e-05"
$ws.Range("L7").Value = 62.02113542173853
$ws.Range("M7").Value = "[34.27120813532028, 89.77106270815679]"
$ws.Range("N7").Value = [double]"4.736578001085512e-05"
$ws.Range("O7").Value = [double]"4.736578001085512e-05"
$ws.Range("P7").Value = 2.836553126664658
$ws.Range("Q7").Value = "[2.3082372449798885, 3.364869008349428]"
$ws.Range("R7").Value = [double]"4.241051954068098e-14"
$ws.Range("S7").Value = [double]"4.241051954068098e-14"
$ws.Range("T7").Value = 71.96295382149904
$ws.Range("U7").Value = "[55.616587610932555, 88.30932003206553]"
$ws.Range("V7").Value = [double]"1.955480222193273e-11"
$ws.Range("W7").Value = [double]"1.955480222193273e-11"
$ws.Range("X7").Value = 12.83603603603616
$ws.Range("Y7").Value = 10.86846846846857
$ws.Range("Z7").Value = 14.80360360360374
$ws.Range("F8").Value = 23.40000000000022
$ws.Range("H8").Value = 0.01247986825177205
$ws.Range("I8").Value = 0.01247986825177205
$ws.Range("L8").Value = 37.91093042011876
$ws.Range("M8").Value = "[6.489391023234674, 69.33246981700285]"
$ws.Range("N8").Value = 0.01914766578739835
$ws.Range("O8").Value = 0.01914766578739835
$ws.Range("P8").Value = 2.849132076228581
$ws.Range("Q8").Value = "[1.7547634641672718, 3.9435006882898906]"
$ws.Range("R8").Value = [double]"4.070964708002833e-06"
$ws.Range("S8").Value = [double]"4.070964708002833e-06"
$ws.Range("T8").Value = 62.64069322605636
$ws.Range("U8").Value = "[45.1171971598301, 80.16418929228261]"
$ws.Range("V8").Value = [double]"5.156276294826512e-09"
$ws.Range("W8").Value = [double]"5.156276294826512e-09"
$ws.Range("X8").Value = 12.78918918918931
$ws.Range("Y8").Value = 8.713513513513593
$ws.Range("Z8").Value = 16.86486486486502
$ws.Range("F9").Value = 22.91000000000014
$ws.Range("H9").Value = [double]"1.375195074260205e-06"
$ws.Range("I9").Value = [double]"1.375195074260205e-06"
$ws.Range("L9").Value = 70.60374902153667
$ws.Range("M9").Value = "[43.44273234789749, 97.76476569517585]"
$ws.Range("N9").Value = [double]"4.182617233894348e-06"
$ws.Range("O9").Value = [double]"4.182617233894348e-06"
$ws.Range("P9").Value = -3.094421592725081
$ws.Range("Q9").Value = "[-3.5598427265902357, -2.6290004588599265]"
$ws.Range("R9").Value = 0
$ws.Range("S9").Value = 0
$ws.Range("T9").Value = 70.7951350183302
$ws.Range("U9").Value = "[53.90611837568764, 87.68415166097276]"
$ws.Range("V9").Value = [double]"7.889244812986362e-11"
$ws.Range("W9").Value = [double]"7.889244812986362e-11"
$ws.Range("X9").Value = 11.28300300300307
$ws.Range("Y9").Value = 9.585965965966023
$ws.Range("Z9").Value = 12.98004004004012
$ws.Range("F10").Value = 22.91000000000014
$ws.Range("H10").Value = 0.000167853920175598
$ws.Range("I10").Value = 0.000167853920175598
$ws.Range("L10").Value = 47.06135118178095
$ws.Range("M10").Value = "[20.325471664704466, 73.79723069885743]"
$ws.Range("N10").Value = 0.0009286932663901126
$ws.Range("O10").Value = 0.0009286932663901126
$ws.Range("P10").Value = -2.830263651882697
$ws.Range("Q10").Value = "[-3.4843690292066984, -2.176158274558695]"
$ws.Range("R10").Value = [double]"3.216449329102034e-11"
$ws.Range("S10").Value = [double]"3.216449329102034e-11"
$ws.Range("T10").Value = 64.5495674868622
$ws.Range("U10").Value = "[49.77115522076834, 79.32797975295605]"
$ws.Range("V10").Value = [double]"2.455191605577056e-11"
$ws.Range("W10").Value = [double]"2.455191605577056e-11"
$ws.Range("X10").Value = 10.31981981981988
$ws.Range("Y10").Value = 7.934794794794842
$ws.Range("Z10").Value = 12.70484484484493
$ws.Range("F11").Value = 22.91000000000014
$ws.Range("H11").Value = 0.002744090536726529
$ws.Range("I11").Value = 0.002744090536726529
$ws.Range("L11").Value = 47.69740062186978
$ws.Range("M11").Value = "[13.295633439051386, 82.09916780468816]"
$ws.Range("N11").Value = 0.007651762604040213
$ws.Range("O11").Value = 0.007651762604040213
$ws.Range("P11").Value = -2.880579450138389
$ws.Range("Q11").Value = "[-3.723369070921236, -2.0377898293555416]"
$ws.Range("R11").Value = [double]"1.514414482706172e-08"
$ws.Range("S11").Value = [double]"1.514414482706172e-08"
$ws.Range("T11").Value = 73.65710104374222
$ws.Range("U11").Value = "[54.86488297009977, 92.44931911738466]"
$ws.Range("V11").Value = [double]"4.913016660168523e-10"
$ws.Range("W11").Value = [double]"4.913016660168523e-10"
$ws.Range("X11").Value = 10.50328328328335
$ws.Range("Y11").Value = 7.430270270270316
$ws.Range("Z11").Value = 13.57629629629638
$ws.Range("F12").Value = 22.91000000000014
$ws.Range("H12").Value = [double]"9.053289895533112e-09"
$ws.Range("I12").Value = [double]"9.053289895533112e-09"
$ws.Range("L12").Value = 73.4478441251764
$ws.Range("M12").Value = "[47.77793413697426, 99.11775411337854]"
$ws.Range("N12").Value = [double]"7.014566874019579e-07"
$ws.Range("O12").Value = [double]"7.014566874019579e-07"
$ws.Range("P12").Value = -2.138421425866927
$ws.Range("Q12").Value = "[-2.503210963220696, -1.7736318885131577]"
$ws.Range("R12").Value = [double]"2.220446049250313e-15"
$ws.Range("S12").Value = [double]"2.220446049250313e-15"
$ws.Range("T12").Value = 63.67749138237657
$ws.Range("U12").Value = "[50.046714217478694, 77.30826854727445]"
$ws.Range("V12").Value = [double]"3.382849556032852e-12"
$ws.Range("W12").Value = [double]"3.382849556032852e-12"
$ws.Range("X12").Value = 7.797197197197248
$ws.Range("Y12").Value = 6.467087087087132
$ws.Range("Z12").Value = 9.127307307307364
